$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-22 14:23:19"
$wsZhCn.Range("G5").Value = "2016-02-22 14:24:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-22 14:23:32"
$wsDeDe.Range("G5").Value = "2016-02-22 14:24:24"
